# Update cryptocurrency price/volume snapshot (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.728.48'
$ws.Range('E2').Value = '  -1.00%  '

$ws.Range('D3').Value = '1.776.07'
$ws.Range('E3').Value = '  -0.96%  '

$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').Value = '''223.95'
$ws.Range('E5').Value = '  +0.69%  '

$ws.Range('E6').Value = '  -0.96%  '

$ws.Range('D8').Value = '''32.16'
$ws.Range('E8').Value = '  +1.15%  '

$ws.Range('E9').Value = '  +1.80%  '

$ws.Range('D10').Value = '''0.0679'
$ws.Range('E10').Value = '  -5.22%  '

$ws.Range('D11').Value = '''0.0935'
$ws.Range('E11').Value = '  +1.41%  '

$ws.Range('D12').Value = '2.031.18'
$ws.Range('E12').Value = '  -0.97%  '

$ws.Range('D13').Value = '''11.22'
$ws.Range('E13').Value = '  +4.37%  '

$ws.Range('D14').Value = '1.759.80'
$ws.Range('E14').Value = '  -1.89%  '

$ws.Range('D15').Value = '33.718.80'
$ws.Range('E15').Value = '  -1.01%  '

$ws.Range('D16').Value = '''0.608'
$ws.Range('E16').Value = '  -3.45%  '

$ws.Range('E17').Value = '  -2.61%  '

$ws.Range('D18').Value = '''66.47'
$ws.Range('E18').Value = '  -2.37%  '

$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '''238.18'
$ws.Range('E19').Value = '  -3.05%  '

$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0772'
$ws.Range('E20').Value = '  -1.45%  '

$ws.Range('E21').Value = '  +0.10%  '

$ws.Range('E22').Value = '  -1.85%  '

$ws.Range('E23').Value = '  -1.92%  '

$ws.Range('D24').Value = '''2.07'
$ws.Range('E24').Value = '  -1.80%  '

$ws.Range('D25').Value = '''159.92'
$ws.Range('E25').Value = '  +0.90%  '

$ws.Range('D26').Value = '''16.07'
$ws.Range('E26').Value = '  -2.23%  '

$ws.Range('E27').Value = '  -0.10%  '

$ws.Range('E28').Value = '  -0.13%  '

$ws.Range('E29').Value = '  +0.15%  '

$ws.Range('E30').Value = '  +1.27%  '

$ws.Range('D31').Value = '''0.0510'
$ws.Range('E31').Value = '  -1.65%  '

$ws.Range('D32').Value = '''3.58'
$ws.Range('E32').Value = '  -3.11%  '

$ws.Range('E33').Value = '  -0.65%  '

$ws.Range('E34').Value = '  -2.01%  '

$ws.Range('D35').Value = '1.382.86'
$ws.Range('E35').Value = '  -1.94%  '

$ws.Range('D36').Value = '''0.645'
$ws.Range('E36').Value = '  +0.21%  '

$ws.Range('E37').Value = '  -2.61%  '

$ws.Range('E38').Value = '  -1.43%  '

$ws.Range('D39').Value = '''2.24'
$ws.Range('E39').Value = '  +5.51%  '

$ws.Range('E40').Value = '  +0.64%  '

$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '''0.907'
$ws.Range('E41').Value = '  -4.04%  '

$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '''78.05'
$ws.Range('E42').Value = '  -2.54%  '

$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '''2.66'
$ws.Range('E43').Value = '  -2.64%  '

$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '''13.50'
$ws.Range('E44').Value = '  +13.97%  '

$ws.Range('E45').Value = '  +3.80%  '

$ws.Range('D46').Value = '''0.0498'
$ws.Range('E46').Value = '  +1.02%  '

$ws.Range('E47').Value = '  +11.93%  '

$ws.Range('D48').Value = '''107.21'
$ws.Range('E48').Value = '  +1.47%  '

$ws.Range('E49').Value = '  -1.77%  '

$ws.Range('D50').Value = '1.930.84'
$ws.Range('E50').Value = '  -0.89%  '

$ws.Range('D51').Value = '''1.00'
$ws.Range('E51').Value = '  +0.20%  '
